$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.972.89"
$ws.Range("D3").Value = "1.645.57"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'214.97"
$ws.Range("E5").Value = "  +2.57%  "
$ws.Range("D6").Value = "'0.5217"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'0.2607"
$ws.Range("D9").Value = "'0.06368"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("D10").Value = "'20.69"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("D11").Value = "'0.07689"
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.650.22"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.426"
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").Value = "1.868.54"
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("D15").Value = "'0.5513"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "0.0₅8261"
$ws.Range("E16").Value = "  +3.42%  "
$ws.Range("D17").Value = "'64.69"
$ws.Range("E17").Value = "  -2.39%  "
$ws.Range("D18").Value = "25.985.87"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "'4.702"
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("E21").Value = "  +1.27%  "
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("D23").Value = "'6.257"
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").Value = "'144.18"
$ws.Range("E25").Value = "  -3.66%  "
$ws.Range("D26").Value = "'0.1245"
$ws.Range("E26").Value = "  +0.88%  "
$ws.Range("D27").Value = "'7.396"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "'15.96"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("D29").Value = "'1.395"
$ws.Range("E29").Value = "  +2.61%  "
$ws.Range("D30").Value = "'0.05910"
$ws.Range("E30").Value = "  -5.62%  "
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("D33").Value = "'3.395"
$ws.Range("E33").Value = "  -2.64%  "
$ws.Range("E34").Value = "  +0.83%  "
$ws.Range("D35").Value = "'0.9922"
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").Value = "'2.748"
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("D38").Value = "'0.5636"
$ws.Range("E38").Value = "  -5.43%  "
$ws.Range("D39").Value = "'0.01604"
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("D40").Value = "'5.874"
$ws.Range("E40").Value = "  -2.89%  "
$ws.Range("D41").Value = "'0.8531"
$ws.Range("E41").Value = "  -0.96%  "
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("D43").Value = "1.035.60"
$ws.Range("E43").Value = "  -6.54%  "
$ws.Range("D44").Value = "'98.73"
$ws.Range("E44").Value = "  -1.91%  "
$ws.Range("D45").Value = "1.792.16"
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈108"
$ws.Range("E46").Value = "  -2.28%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'55.66"
$ws.Range("E47").Value = "  +0.73%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.004"
$ws.Range("E48").Value = "  +0.27%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.039"
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.05145"
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.4216"
$ws.Range("E51").Value = "  -0.49%  "
